$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.26%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.00%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.121"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.19%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'-1.40%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.61%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").Value = "'2.488"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.67%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9021"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.41%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1110"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'9.78%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1765"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.79%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09192"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.08%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04201"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.33%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.26%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001250"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.79%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005873"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.21%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.353"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.02%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.253"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.19%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-0.95%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.548"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-7.07%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1360"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.64%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2682"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-10.46%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04072"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.26%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001231"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.67%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004094"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.61%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.07%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D38").Value = "'0.02402"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'1.77%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05189"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.46%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007783"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.71%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1302"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.88%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006748"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.84%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-0.64%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008764"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.23%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3333"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.54%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007012"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'7.52%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.03079"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1,314.09%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004202"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-40.03%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.07%"
$ws.Range("E51").Style = "Normal"
